$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

$emptyCols = @(2, 9, 10, 11, 12, 13, 17, 19, 23)

$ws.Cells.Item($row, 1).Value = "2024-09-01 21:15:40"
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 1
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 14).Value = 10
$ws.Cells.Item($row, 15).Value = 6
$ws.Cells.Item($row, 16).Value = 3
$ws.Cells.Item($row, 18).Value = 5
$ws.Cells.Item($row, 20).Value = 50
$ws.Cells.Item($row, 21).Value = 0.3333333333333333
$ws.Cells.Item($row, 22).Value = "D:\Repositorio\jonatha1992\Predictor_ruleta_ejecutable\Data\Electromecanica.xlsx"
$ws.Cells.Item($row, 24).Value = "No es Simulación"
$ws.Cells.Item($row, 25).Value = 3

foreach ($col in $emptyCols) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'"
    $cell.Style = "Normal"
}
